$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65, shifting existing rows 65-98 down to 66-99
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with this week's data point
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C65").Value = 'Arica y Parinacota'
$ws.Range("D65").Value = 45029
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112012
$ws.Range("G65").Value = 'Espinaca'
$ws.Range("H65").Value = 'Sin especificar'
$ws.Range("I65").Value = 'Primera'
$ws.Range("J65").Value = 200
$ws.Range("K65").Value = 4000
$ws.Range("L65").Value = 4500
$ws.Range("M65").Value = 4250
$ws.Range("N65").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O65").Value = 'Región de Arica y Parinacota'
$ws.Range("P65").Value = 1417
$ws.Range("Q65").Value = 3
$ws.Range("R65").Value = 'Hortaliza'
